$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.217.56"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").Value = "1.660.64"
$ws.Range("E3").Value = "  -0.65%  "

$__style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $__style
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  -1.06%  "

$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5200"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  -1.36%  "

$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  +0.02%  "

$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2643"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  -1.06%  "

$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06276"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -1.61%  "

$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.83"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  -3.99%  "

$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07774"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  -0.44%  "

$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.482"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "1.636.49"
$ws.Range("E13").Value = "  -2.20%  "

$ws.Range("D14").Value = "1.886.88"
$ws.Range("E14").Value = "  -0.69%  "

$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5469"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "0.0₅8148"
$ws.Range("E16").Value = "  -1.75%  "

$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.02"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").Value = "26.210.35"

$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  +0.08%  "

$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.612"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  -2.70%  "

$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.90"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  -0.55%  "

$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  -2.60%  "

$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.010"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  -4.21%  "

$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +0.00%  "

$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.08"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +0.65%  "

$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1225"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  -2.76%  "

$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.293"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  -1.33%  "

$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.14"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  -0.51%  "

$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.436"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  +1.25%  "

$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05960"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  -3.25%  "

$ws.Range("E31").Value = "  -0.86%  "

$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.550"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  -1.69%  "

$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.276"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  -3.44%  "

$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.585"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  -5.56%  "

$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9610"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  -4.18%  "

$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.417"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  -0.12%  "

$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.765"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  -0.11%  "

$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5699"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -5.53%  "

$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01593"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  -1.01%  "

$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.987"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  -0.46%  "

$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8496"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  -0.62%  "

$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.005.05"
$ws.Range("E43").Value = "  -8.33%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.45"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "1.800.93"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("E46").Value = "  +9.12%  "

$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.58"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  -2.52%  "

$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +0.38%  "

$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.036"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  -1.34%  "

$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4335"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +2.45%  "

$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05154"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  -0.95%  "
